$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 499.83334
$ws.Range("I2").Value = 699.75
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 699.75
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -586.75
$ws.Range("N2").Value = -326

$ws.Range("H4").Value = 148.25
$ws.Range("I4").Value = 148.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 148.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -34.25

$ws.Range("H11").Value = 80.76922999999999
$ws.Range("I11").Value = 80.76922999999999
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 80.76922999999999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 59.23077000000001

$ws.Range("H12").Value = 450
$ws.Range("I12").Value = 201
$ws.Range("J12").Value = 699
$ws.Range("K12").Value = 201
$ws.Range("L12").Value = 699
$ws.Range("M12").Value = -31
$ws.Range("N12").Value = -1039

$ws.Range("H17").Value = 53716.117
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 56998.375
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 170995.125
$ws.Range("M17").Value = -3432
$ws.Range("N17").Value = -171331.125

$ws.Range("H28").Value = 276.29413
$ws.Range("I28").Value = 282.3125
$ws.Range("J28").Value = 180
$ws.Range("K28").Value = 282.3125
$ws.Range("L28").Value = 180
$ws.Range("M28").Value = 202.6875
$ws.Range("N28").Value = -1150

$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 1000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -771

$ws.Range("H70").Value = 1064.8
$ws.Range("I70").Value = 899.3333
$ws.Range("J70").Value = 1135.7142
$ws.Range("K70").Value = 2697.9999
$ws.Range("L70").Value = 3407.1426
$ws.Range("M70").Value = -2427.9999
$ws.Range("N70").Value = -3947.1426

$ws.Range("H73").Value = 1064.8
$ws.Range("I73").Value = 899.3333
$ws.Range("J73").Value = 1135.7142
$ws.Range("K73").Value = 2697.9999
$ws.Range("L73").Value = 3407.1426
$ws.Range("M73").Value = -1761.9999
$ws.Range("N73").Value = -5279.142599999999

$ws.Range("H80").Value = 527.8
$ws.Range("I80").Value = 450
$ws.Range("J80").Value = 579.6667
$ws.Range("K80").Value = 1350
$ws.Range("L80").Value = 1739.0001
$ws.Range("M80").Value = -352
$ws.Range("N80").Value = -3735.0001

$ws.Range("H83").Value = 527.8
$ws.Range("I83").Value = 450
$ws.Range("J83").Value = 579.6667
$ws.Range("K83").Value = 4050
$ws.Range("L83").Value = 5217.0003
$ws.Range("M83").Value = 942
$ws.Range("N83").Value = -15201.0003

$ws.Range("H97").Value = 2250
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2250
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6750
$ws.Range("N97").Value = -7742

$ws.Range("H100").Value = 3301.4375
$ws.Range("I100").Value = 2856.7778
$ws.Range("J100").Value = 3873.1428
$ws.Range("K100").Value = 2856.7778
$ws.Range("L100").Value = 3873.1428
$ws.Range("M100").Value = -2315.7778
$ws.Range("N100").Value = -4955.1428

$ws.Range("H112").Value = 3365.923
$ws.Range("I112").Value = 2319.6
$ws.Range("J112").Value = 3615.0476
$ws.Range("K112").Value = 6958.799999999999
$ws.Range("L112").Value = 10845.1428
$ws.Range("M112").Value = -5850.799999999999
$ws.Range("N112").Value = -13061.1428

$ws.Range("H113").Value = 4505.9287
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 4737.154
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 4737.154
$ws.Range("M113").Value = 1754
$ws.Range("N113").Value = -11245.154

$ws.Range("H125").Value = 27997.416
$ws.Range("I125").Value = 36814.11
$ws.Range("J125").Value = 1547.3334
$ws.Range("K125").Value = 331326.99
$ws.Range("L125").Value = 13926.0006
$ws.Range("M125").Value = -328866.99
$ws.Range("N125").Value = -18846.0006

$ws.Range("H132").Value = 3387.7856
$ws.Range("I132").Value = 1343.9667
$ws.Range("J132").Value = 8497.333000000001
$ws.Range("K132").Value = 4031.9001
$ws.Range("L132").Value = 25491.999
$ws.Range("M132").Value = -1501.9001
$ws.Range("N132").Value = -30551.999

$ws.Range("H138").Value = 1853.738
$ws.Range("I138").Value = 1004.8889
$ws.Range("J138").Value = 3381.6667
$ws.Range("K138").Value = 3014.6667
$ws.Range("L138").Value = 10145.0001
$ws.Range("M138").Value = 2125.3333
$ws.Range("N138").Value = -20425.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 203.23077
$ws.Range("I5").Value = 178
$ws.Range("J5").Value = 232.66667
$ws.Range("K5").Value = 178
$ws.Range("L5").Value = 232.66667
$ws.Range("M5").Value = -66
$ws.Range("N5").Value = -456.66667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 203.23077
$ws.Range("I4").Value = 178
$ws.Range("J4").Value = 232.66667
$ws.Range("K4").Value = 178
$ws.Range("L4").Value = 232.66667
$ws.Range("M4").Value = -63
$ws.Range("N4").Value = -462.66667

$ws.Range("H134").Value = 1399
$ws.Range("I134").Value = 1399
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4197
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1975.6666
$ws.Range("I105").Value = 1978.3334
$ws.Range("J105").Value = 1973
$ws.Range("K105").Value = 1978.3334
$ws.Range("L105").Value = 1973
$ws.Range("M105").Value = -231.3334
$ws.Range("N105").Value = -5467

$ws.Range("H132").Value = 3023.9355
$ws.Range("I132").Value = 3139.2593
$ws.Range("J132").Value = 2245.5
$ws.Range("K132").Value = 9417.777900000001
$ws.Range("L132").Value = 6736.5
$ws.Range("M132").Value = -6887.777900000001
$ws.Range("N132").Value = -11796.5

$ws.Range("H134").Value = 2444.682
$ws.Range("I134").Value = 2556.625
$ws.Range("J134").Value = 2146.1667
$ws.Range("K134").Value = 7669.875
$ws.Range("L134").Value = 6438.500100000001
$ws.Range("M134").Value = -5134.875
$ws.Range("N134").Value = -11508.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 143185.28
$ws.Range("I11").Value = 83
$ws.Range("J11").Value = 250512
$ws.Range("K11").Value = 249
$ws.Range("L11").Value = 751536
$ws.Range("M11").Value = -109
$ws.Range("N11").Value = -751816

$ws.Range("H19").Value = 150
$ws.Range("I19").Value = 150
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 450
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -276
$ws.Range("N19").ClearContents()

$ws.Range("H52").Value = 1114.3334
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1114.3334
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 3343.0002
$ws.Range("N52").Value = -3875.0002

$ws.Range("H76").Value = 14989.8
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 14989.8
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 44969.39999999999
$ws.Range("N76").Value = -45735.39999999999

$ws.Range("H79").Value = 14989.8
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 14989.8
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 44969.39999999999
$ws.Range("N79").Value = -47621.39999999999

$ws.Range("H133").Value = 12108.777
$ws.Range("I133").Value = 10996
$ws.Range("J133").Value = 12999
$ws.Range("K133").Value = 32988
$ws.Range("L133").Value = 38997
$ws.Range("M133").Value = -27928
$ws.Range("N133").Value = -49117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 23905.268
$ws.Range("I61").Value = 1546.909
$ws.Range("J61").Value = 85390.75
$ws.Range("K61").Value = 1546.909
$ws.Range("L61").Value = 85390.75
$ws.Range("M61").Value = -1344.909
$ws.Range("N61").Value = -85794.75

$ws.Range("H100").Value = 26992.75
$ws.Range("I100").Value = 7654.5557
$ws.Range("J100").Value = 51856.145
$ws.Range("K100").Value = 7654.5557
$ws.Range("L100").Value = 51856.145
$ws.Range("M100").Value = -7113.5557
$ws.Range("N100").Value = -52938.145

$ws.Range("H113").Value = 23905.268
$ws.Range("I113").Value = 1546.909
$ws.Range("J113").Value = 85390.75
$ws.Range("K113").Value = 1546.909
$ws.Range("L113").Value = 85390.75
$ws.Range("M113").Value = 623.0909999999999
$ws.Range("N113").Value = -89730.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 24780.666
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 24780.666
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 24780.666
$ws.Range("N64").Value = -25276.666

$ws.Range("H67").Value = 24780.666
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 24780.666
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 24780.666
$ws.Range("N67").Value = -26496.666

$ws.Range("H96").Value = 3589.2856
$ws.Range("I96").Value = 3711.5386
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 3711.5386
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -2338.5386
$ws.Range("N96").Value = -4746

$ws.Range("H132").Value = 15802.6875
$ws.Range("I132").Value = 14446.392
$ws.Range("J132").Value = 46997.5
$ws.Range("K132").Value = 43339.176
$ws.Range("L132").Value = 140992.5
$ws.Range("M132").Value = -40809.176
$ws.Range("N132").Value = -146052.5
